$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Introduction" row (row 3), shifting rows 4-7 up
$ws.Rows("3:3").Delete()

# Update Hours values for the shifted rows
$ws.Range("B3").Value = 18
$ws.Range("B4").Value = 0.5
$ws.Range("B5").Value = 6.5
$ws.Range("B6").Value = 10

# Update selection
$ws.Range("C7").Select()
